$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-3.10%"
$ws.Range("D3").Value = "'54.32"
$ws.Range("E3").Value = "'10.13%"
$ws.Range("D4").Value = "'5.105"
$ws.Range("E4").Value = "'-3.74%"
$ws.Range("D5").Value = "'0.07932"
$ws.Range("E5").Value = "'-1.80%"
$ws.Range("D6").Value = "'4.574"
$ws.Range("E6").Value = "'-0.30%"
$ws.Range("D7").Value = "'1.401"
$ws.Range("E7").Value = "'4.19%"
$ws.Range("D8").Value = "'1.666"
$ws.Range("E8").Value = "'0.86%"
$ws.Range("D9").Value = "'0.1237"
$ws.Range("E9").Value = "'-2.42%"
$ws.Range("D10").Value = "'0.2002"
$ws.Range("E10").Value = "'1.44%"
$ws.Range("D11").Value = "'0.04731"
$ws.Range("E11").Value = "'1.07%"
$ws.Range("D12").Value = "'0.09469"
$ws.Range("E12").Value = "'-2.30%"
$ws.Range("D13").Value = "'0.1045"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("D14").Value = "'0.001273"
$ws.Range("E14").Value = "'-3.92%"
$ws.Range("D15").Value = "'0.005850"
$ws.Range("E15").Value = "'-0.59%"
$ws.Range("D16").Value = "'3.339"
$ws.Range("E16").Value = "'-0.31%"
$ws.Range("E17").Value = "'-2.47%"
$ws.Range("D18").Value = "'0.3420"
$ws.Range("E18").Value = "'-2.90%"
$ws.Range("D19").Value = "'8.396"
$ws.Range("E19").Value = "'3.11%"
$ws.Range("D20").Value = "'0.1358"
$ws.Range("E20").Value = "'-1.65%"
$ws.Range("E21").Value = "'-5.51%"
$ws.Range("D22").Value = "'0.04184"
$ws.Range("E22").Value = "'-0.50%"
$ws.Range("D23").Value = "'0.001257"
$ws.Range("E23").Value = "'-3.02%"
$ws.Range("D24").Value = "'0.004094"
$ws.Range("E24").Value = "'-5.28%"
$ws.Range("D25").Value = "'0.0001347"
$ws.Range("E25").Value = "'-0.09%"
$ws.Range("D26").Value = "'0.0003538"
$ws.Range("E26").Value = "'0.14%"
$ws.Range("D38").Value = "'0.02625"
$ws.Range("E38").Value = "'-3.46%"
$ws.Range("D39").Value = "'0.05950"
$ws.Range("E39").Value = "'-0.38%"
$ws.Range("D40").Value = "'0.01083"
$ws.Range("E40").Value = "'0.48%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007978"
$ws.Range("E41").Value = "'-0.70%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1661"
$ws.Range("E42").Value = "'13.11%"
$ws.Range("D43").Value = "'0.008195"
$ws.Range("E43").Value = "'3.73%"
$ws.Range("D44").Value = "'0.008299"
$ws.Range("E44").Value = "'5.37%"
$ws.Range("D45").Value = "'0.3438"
$ws.Range("E45").Value = "'6.98%"
$ws.Range("D46").Value = "'0.00007332"
$ws.Range("E46").Value = "'5.33%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("D48").Value = "'0.05550"
$ws.Range("E48").Value = "'-0.78%"
$ws.Range("D49").Value = "'0.002619"
$ws.Range("E49").Value = "'-34.42%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'0.14%"
